$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.478.62'
$ws.Range("E2").Value = '  +5.49%  '

# Row 3
$ws.Range("D3").Value = '3.624.19'
$ws.Range("E3").Value = '  +5.32%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.44'
$ws.Range("E5").Value = '  +2.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.50'
$ws.Range("E6").Value = '  +3.70%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.645'
$ws.Range("E7").Value = '  +2.04%  '

# Row 8
$ws.Range("D8").Value = '3.614.64'
$ws.Range("E8").Value = '  +5.22%  '

# Row 9
$ws.Range("E9").Value = '  -0.03%  '

# Row 10
$ws.Range("E10").Value = '  +2.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.666'
$ws.Range("E11").Value = '  +3.11%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.46'
$ws.Range("E12").Value = '  +4.04%  '

# Row 13
$ws.Range("E13").Value = '  +4.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.90'
$ws.Range("E14").Value = '  +5.03%  '

# Row 15
$ws.Range("D15").Value = '4.194.03'
$ws.Range("E15").Value = '  +5.09%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.68'
$ws.Range("E16").Value = '  +5.38%  '

# Row 17
$ws.Range("D17").Value = '3.612.51'
$ws.Range("E17").Value = '  +4.99%  '

# Row 18
$ws.Range("D18").Value = '70.352.20'
$ws.Range("E18").Value = '  +5.46%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.64'
$ws.Range("E19").Value = '  +4.67%  '

# Row 20
$ws.Range("E20").Value = '  +0.76%  '

# Row 21
$ws.Range("E21").Value = '  +4.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '491.49'
$ws.Range("E22").Value = '  +1.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.54'
$ws.Range("E23").Value = '  +15.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.34'
$ws.Range("E24").Value = '  +5.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.46'
$ws.Range("E25").Value = '  +2.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.04'
$ws.Range("E26").Value = '  +1.55%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.13'
$ws.Range("E27").Value = '  +6.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.34'
$ws.Range("E28").Value = '  +2.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.64'
$ws.Range("E29").Value = '  +4.78%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.16'
$ws.Range("E30").Value = '  +5.74%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.75'
$ws.Range("E31").Value = '  +9.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '639.57'
$ws.Range("E32").Value = '  +7.84%  '

# Row 33
$ws.Range("E33").Value = '  +5.17%  '

# Row 34
$ws.Range("E34").Value = '  +7.40%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.86'
$ws.Range("E35").Value = '  +2.08%  '

# Row 36
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.78'
$ws.Range("E36").Value = '  +6.52%  '

# Row 37
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0822'
$ws.Range("E37").Value = '  +6.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.410'
$ws.Range("E38").Value = '  +5.97%  '

# Row 39
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  -0.38%  '

# Row 42
$ws.Range("D42").Value = '3.302.34'
$ws.Range("E42").Value = '  +3.23%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.12'
$ws.Range("E43").Value = '  +6.63%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  +9.85%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0452'
$ws.Range("E45").Value = '  +5.16%  '

# Row 46
$ws.Range("E46").Value = '  +2.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("E47").Value = '  +2.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.12'
$ws.Range("E48").Value = '  +4.68%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.70'
$ws.Range("E49").Value = '  -3.25%  '

# Row 50
$ws.Range("E50").Value = '  +4.14%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.06%  '
